# Update Name of Algo
# Apply updated imputed values (RandomForest re-run) to columns B and D
# for the specific rows that changed, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = 6.026799999999998
$ws.Range("D11").Value = -8.212899999999992
$ws.Range("B12").Value = 4.949799999999998
$ws.Range("B15").Value = 5.189699999999997
$ws.Range("D23").Value = -8.098000000000004
$ws.Range("B27").Value = 6.385400000000007
$ws.Range("B28").Value = 5.886400000000005
$ws.Range("D28").Value = -8.530999999999995
$ws.Range("B31").Value = 5.094300000000001
$ws.Range("B32").Value = 6.459099999999995
$ws.Range("D32").Value = -8.509599999999985
$ws.Range("D34").Value = -8.119200000000001
$ws.Range("B36").Value = 10.0936
$ws.Range("D36").Value = -7.144600000000002
$ws.Range("D37").Value = -8.085600000000005
$ws.Range("B38").Value = 4.818699999999996
$ws.Range("D42").Value = -9.148499999999995
$ws.Range("B46").Value = 6.173800000000004
$ws.Range("D49").Value = -8.023899999999996
$ws.Range("B54").Value = 5.056700000000004
$ws.Range("D54").Value = -8.140299999999995
$ws.Range("B55").Value = 5.923299999999998
$ws.Range("B56").Value = 4.458899999999996
$ws.Range("B67").Value = 5.216599999999996
$ws.Range("B69").Value = 5.342599999999996
$ws.Range("B72").Value = 5.275000000000005
$ws.Range("B73").Value = 8.680900000000001
$ws.Range("D78").Value = -7.478700000000003
$ws.Range("D80").Value = -7.798800000000004
$ws.Range("B83").Value = 5.633599999999998
$ws.Range("B86").Value = 5.029500000000002
$ws.Range("B91").Value = 5.104499999999995
$ws.Range("B93").Value = 5.179999999999999
$ws.Range("D97").Value = -8.462099999999994
$ws.Range("B99").Value = 4.798500000000002
$ws.Range("D99").Value = -8.431399999999998
$ws.Range("D100").Value = -8.202400000000001
$ws.Range("D101").Value = -7.707000000000003
$ws.Range("B104").Value = 9.916900000000004
$ws.Range("B105").Value = 8.492000000000003
